$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 325 (shifts existing rows 325:380 down to 326:381)
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row 325 with the new record
$ws.Cells.Item(325, 1).Value2 = 5
$ws.Cells.Item(325, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(325, 3).Value2 = "Maule"
$ws.Cells.Item(325, 4).Value2 = 45209
$ws.Cells.Item(325, 5).Value2 = 7
$ws.Cells.Item(325, 6).Value2 = 100112017
$ws.Cells.Item(325, 7).Value2 = "Apio"
$ws.Cells.Item(325, 8).Value2 = "Americana (o)"
$ws.Cells.Item(325, 9).Value2 = "Primera"
$ws.Cells.Item(325, 10).Value2 = 700
$ws.Cells.Item(325, 11).Value2 = 6000
$ws.Cells.Item(325, 12).Value2 = 6000
$ws.Cells.Item(325, 13).Value2 = 6000
$ws.Cells.Item(325, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(325, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(325, 16).Value2 = 1000
$ws.Cells.Item(325, 17).Value2 = 6
$ws.Cells.Item(325, 18).Value2 = "Hortaliza"
